# Auto-generated edit script: rebuild the Export sheet data to match target state
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure account-number column (A) for the data rows is formatted as Text so that
# values with leading zeros (e.g. "005135532") are preserved as text, not numbers.
$ws.Range("A2:A248").NumberFormat = "@"

$rowCount = 247
$data = New-Object 'object[,]' $rowCount,3
$data[0,0] = '005135532'
$data[0,1] = 'FELIPE'
$data[0,2] = 190000
$data[1,0] = '004211368'
$data[1,1] = 'ILTON'
$data[1,2] = 139849.01
$data[2,0] = '004431546'
$data[2,1] = 'GABRIELA'
$data[2,2] = 100000
$data[3,0] = '004239624'
$data[3,1] = 'NINA'
$data[3,2] = 13190.59
$data[4,0] = '005142661'
$data[4,1] = 'SABRINA'
$data[4,2] = 11000
$data[5,0] = '001879977'
$data[5,1] = 'THAISSA'
$data[5,2] = 10992.16
$data[6,0] = '005624274'
$data[6,1] = 'CLAYTON'
$data[6,2] = 4297.96
$data[7,0] = '004313254'
$data[7,1] = 'GUSTAVO'
$data[7,2] = 4292
$data[8,0] = '004886366'
$data[8,1] = 'RENATO'
$data[8,2] = 3297.65
$data[9,0] = '004368468'
$data[9,1] = 'AHMAD'
$data[9,2] = 2766.45
$data[10,0] = '004213139'
$data[10,1] = 'LEONARDO'
$data[10,2] = 2609.88
$data[11,0] = '005000460'
$data[11,1] = 'MARIANA'
$data[11,2] = 2208.06
$data[12,0] = '002823185'
$data[12,1] = 'SIMONE'
$data[12,2] = 2165.46
$data[13,0] = '001761119'
$data[13,1] = 'BLUEMETRIX'
$data[13,2] = 979.41
$data[14,0] = '004329030'
$data[14,1] = 'DANIELA'
$data[14,2] = 940.23
$data[15,0] = '004392159'
$data[15,1] = 'RODRIGO'
$data[15,2] = 900.21
$data[16,0] = '005685353'
$data[16,1] = 'CARLOS'
$data[16,2] = 767.05
$data[17,0] = '005696595'
$data[17,1] = 'CLUBE'
$data[17,2] = 752.05
$data[18,0] = '004855960'
$data[18,1] = 'CLERIA'
$data[18,2] = 556.35
$data[19,0] = '004220849'
$data[19,1] = 'DULCE'
$data[19,2] = 503.59
$data[20,0] = '008002502'
$data[20,1] = 'JORGEANA'
$data[20,2] = 500
$data[21,0] = '000806386'
$data[21,1] = 'FERNANDA'
$data[21,2] = 492.08
$data[22,0] = '005002457'
$data[22,1] = 'ROSANGELA'
$data[22,2] = 484.08
$data[23,0] = '004488571'
$data[23,1] = 'CARLOS'
$data[23,2] = 454.3
$data[24,0] = '004432579'
$data[24,1] = 'ANA'
$data[24,2] = 446.18
$data[25,0] = '005186167'
$data[25,1] = 'ANDREA'
$data[25,2] = 441.75
$data[26,0] = '004377713'
$data[26,1] = 'DANIELI'
$data[26,2] = 396.98
$data[27,0] = '004508516'
$data[27,1] = 'EDUARDO'
$data[27,2] = 364.49
$data[28,0] = '004355790'
$data[28,1] = 'MINEIA'
$data[28,2] = 323.87
$data[29,0] = '005040864'
$data[29,1] = 'ANDRE'
$data[29,2] = 279.96
$data[30,0] = '004374891'
$data[30,1] = 'RODRIGO'
$data[30,2] = 273.15
$data[31,0] = '004363260'
$data[31,1] = 'LARISSA'
$data[31,2] = 257.51
$data[32,0] = '004515341'
$data[32,1] = 'BRUNO'
$data[32,2] = 235.12
$data[33,0] = '003249855'
$data[33,1] = 'MARINA'
$data[33,2] = 223.31
$data[34,0] = '004472431'
$data[34,1] = 'LUIS'
$data[34,2] = 197.96
$data[35,0] = '003301389'
$data[35,1] = 'EDMUNDO'
$data[35,2] = 191.02
$data[36,0] = '004487016'
$data[36,1] = 'ROGERIO'
$data[36,2] = 176.96
$data[37,0] = '005591536'
$data[37,1] = 'GUSTAVO'
$data[37,2] = 170.31
$data[38,0] = '004508526'
$data[38,1] = 'CASSIO'
$data[38,2] = 153.62
$data[39,0] = '004381415'
$data[39,1] = 'JOAO'
$data[39,2] = 127.55
$data[40,0] = '008070544'
$data[40,1] = 'MARINA'
$data[40,2] = 122.74
$data[41,0] = '008090243'
$data[41,1] = 'GABRIEL'
$data[41,2] = 122.74
$data[42,0] = '004272426'
$data[42,1] = 'RODRIGO'
$data[42,2] = 115.02
$data[43,0] = '005135105'
$data[43,1] = 'BRENNER'
$data[43,2] = 111.75
$data[44,0] = '005646524'
$data[44,1] = 'EVANGELINA'
$data[44,2] = 109.74
$data[45,0] = '004207184'
$data[45,1] = 'CRISTINA'
$data[45,2] = 100.15
$data[46,0] = '008071998'
$data[46,1] = 'ISADORA'
$data[46,2] = 100
$data[47,0] = '005295509'
$data[47,1] = 'BHRUNA'
$data[47,2] = 99.54
$data[48,0] = '005701765'
$data[48,1] = 'F'
$data[48,2] = 98.96
$data[49,0] = '001368670'
$data[49,1] = 'THIAGO'
$data[49,2] = 97.2
$data[50,0] = '004239387'
$data[50,1] = 'LUIZ'
$data[50,2] = 94.87
$data[51,0] = '004431591'
$data[51,1] = 'MARIO'
$data[51,2] = 94.24
$data[52,0] = '004383268'
$data[52,1] = 'LAURA'
$data[52,2] = 93.29
$data[53,0] = '004384258'
$data[53,1] = 'PAULA'
$data[53,2] = 93.29
$data[54,0] = '004466350'
$data[54,1] = 'RAQUEL'
$data[54,2] = 93.1
$data[55,0] = '004536602'
$data[55,1] = 'TATIANY'
$data[55,2] = 92.21
$data[56,0] = '003115072'
$data[56,1] = 'VICTOR'
$data[56,2] = 89.47
$data[57,0] = '004462930'
$data[57,1] = 'WALTER'
$data[57,2] = 87.61
$data[58,0] = '004212132'
$data[58,1] = 'JOAO'
$data[58,2] = 86.38
$data[59,0] = '004809902'
$data[59,1] = 'PEDRO'
$data[59,2] = 85.9
$data[60,0] = '004517080'
$data[60,1] = 'TATIANA'
$data[60,2] = 85.59
$data[61,0] = '004261201'
$data[61,1] = 'ANA'
$data[61,2] = 83.09
$data[62,0] = '004424761'
$data[62,1] = 'PEDRO'
$data[62,2] = 80
$data[63,0] = '004384167'
$data[63,1] = 'DOUGLAS'
$data[63,2] = 79.87
$data[64,0] = '008069255'
$data[64,1] = 'ANGELA'
$data[64,2] = 77.23
$data[65,0] = '000827730'
$data[65,1] = 'LUCIANA'
$data[65,2] = 76.01
$data[66,0] = '005142611'
$data[66,1] = 'GUILHERME'
$data[66,2] = 74.22
$data[67,0] = '001719494'
$data[67,1] = 'LUIS'
$data[67,2] = 73.48
$data[68,0] = '004912314'
$data[68,1] = 'FABRICIO'
$data[68,2] = 71.85
$data[69,0] = '005880251'
$data[69,1] = 'LUIZ'
$data[69,2] = 70.94
$data[70,0] = '004563252'
$data[70,1] = 'FERNANDO'
$data[70,2] = 70.58
$data[71,0] = '004477812'
$data[71,1] = 'DIEGO'
$data[71,2] = 70.23
$data[72,0] = '004877741'
$data[72,1] = 'LUIZ'
$data[72,2] = 70.02
$data[73,0] = '004482090'
$data[73,1] = 'CEZAR'
$data[73,2] = 69.96
$data[74,0] = '005173958'
$data[74,1] = 'VENIA'
$data[74,2] = 68.22
$data[75,0] = '004452507'
$data[75,1] = 'DANIELA'
$data[75,2] = 67.76
$data[76,0] = '004212409'
$data[76,1] = 'RAFAEL'
$data[76,2] = 67.39
$data[77,0] = '004381194'
$data[77,1] = 'ALINNE'
$data[77,2] = 67.06
$data[78,0] = '008032257'
$data[78,1] = 'SARA'
$data[78,2] = 67.03
$data[79,0] = '004457389'
$data[79,1] = 'RAFAEL'
$data[79,2] = 66.93
$data[80,0] = '004749928'
$data[80,1] = 'NADY'
$data[80,2] = 66.86
$data[81,0] = '004242237'
$data[81,1] = 'MARIAH'
$data[81,2] = 66.47
$data[82,0] = '004027477'
$data[82,1] = 'GABRIELA'
$data[82,2] = 64.77
$data[83,0] = '004335251'
$data[83,1] = 'EDMUNDO'
$data[83,2] = 62.82
$data[84,0] = '005924958'
$data[84,1] = 'TIAGO'
$data[84,2] = 62.24
$data[85,0] = '003836362'
$data[85,1] = 'ISABELLA'
$data[85,2] = 61.49
$data[86,0] = '005092207'
$data[86,1] = 'BRUNO'
$data[86,2] = 61.37
$data[87,0] = '005068961'
$data[87,1] = 'JORGE'
$data[87,2] = 61.09
$data[88,0] = '005141215'
$data[88,1] = 'KARINA'
$data[88,2] = 60.56
$data[89,0] = '005890232'
$data[89,1] = 'TAYLA'
$data[89,2] = 59.71
$data[90,0] = '005558076'
$data[90,1] = 'ALEXANDRE'
$data[90,2] = 59.47
$data[91,0] = '004384131'
$data[91,1] = 'ANDRE'
$data[91,2] = 58.71
$data[92,0] = '004382374'
$data[92,1] = 'THEOMAR'
$data[92,2] = 57.75
$data[93,0] = '004974089'
$data[93,1] = 'CELIA'
$data[93,2] = 57.68
$data[94,0] = '004335144'
$data[94,1] = 'EDMUNDO'
$data[94,2] = 57.28
$data[95,0] = '004334062'
$data[95,1] = 'MERG'
$data[95,2] = 57.01
$data[96,0] = '004546050'
$data[96,1] = 'LUIS'
$data[96,2] = 56.88
$data[97,0] = '005079458'
$data[97,1] = 'JONAS'
$data[97,2] = 56.38
$data[98,0] = '004218542'
$data[98,1] = 'JOSE'
$data[98,2] = 56.33
$data[99,0] = '004813166'
$data[99,1] = 'VENIA'
$data[99,2] = 55.69
$data[100,0] = '004552021'
$data[100,1] = 'MARIA'
$data[100,2] = 53.8
$data[101,0] = '004504449'
$data[101,1] = 'KELMA'
$data[101,2] = 53.39
$data[102,0] = '004643880'
$data[102,1] = 'GABRIEL'
$data[102,2] = 53.11
$data[103,0] = '004400640'
$data[103,1] = 'FELIPE'
$data[103,2] = 51.44
$data[104,0] = '004426743'
$data[104,1] = 'GABRIELLE'
$data[104,2] = 51.09
$data[105,0] = '008032413'
$data[105,1] = 'VICTOR'
$data[105,2] = 51.06
$data[106,0] = '004332207'
$data[106,1] = 'IRACY'
$data[106,2] = 48.95
$data[107,0] = '005514036'
$data[107,1] = 'ANA'
$data[107,2] = 48.75
$data[108,0] = '005076418'
$data[108,1] = 'LEONARDO'
$data[108,2] = 48.27
$data[109,0] = '004693308'
$data[109,1] = 'LAURA'
$data[109,2] = 48.25
$data[110,0] = '005152037'
$data[110,1] = 'RODRIGO'
$data[110,2] = 47.9
$data[111,0] = '004260002'
$data[111,1] = 'ERICA'
$data[111,2] = 46.9
$data[112,0] = '004508159'
$data[112,1] = 'FELIPE'
$data[112,2] = 46.28
$data[113,0] = '001731007'
$data[113,1] = 'GUILHERME'
$data[113,2] = 44.73
$data[114,0] = '005981575'
$data[114,1] = 'GLAUCIANE'
$data[114,2] = 44.17
$data[115,0] = '004278212'
$data[115,1] = 'LEONARDO'
$data[115,2] = 43.87
$data[116,0] = '004381095'
$data[116,1] = 'THIAGO'
$data[116,2] = 42.94
$data[117,0] = '004332103'
$data[117,1] = 'JOSE'
$data[117,2] = 42.34
$data[118,0] = '004425965'
$data[118,1] = 'CAROLLINA'
$data[118,2] = 41.94
$data[119,0] = '004290978'
$data[119,1] = 'LARISSA'
$data[119,2] = 40.34
$data[120,0] = '004971448'
$data[120,1] = 'CLOVIS'
$data[120,2] = 40.07
$data[121,0] = '005245032'
$data[121,1] = 'ROSA'
$data[121,2] = 39.91
$data[122,0] = '004238164'
$data[122,1] = 'DANIELA'
$data[122,2] = 39.09
$data[123,0] = '004752519'
$data[123,1] = 'MARCUS'
$data[123,2] = 37.58
$data[124,0] = '004752615'
$data[124,1] = 'LUZIMAR'
$data[124,2] = 37.11
$data[125,0] = '004413537'
$data[125,1] = 'CLAUDIA'
$data[125,2] = 36.55
$data[126,0] = '005077648'
$data[126,1] = 'DUNAS'
$data[126,2] = 36.46
$data[127,0] = '004806286'
$data[127,1] = 'VERA'
$data[127,2] = 35.77
$data[128,0] = '005266369'
$data[128,1] = 'EG'
$data[128,2] = 35.27
$data[129,0] = '004240292'
$data[129,1] = 'MARCO'
$data[129,2] = 34.71
$data[130,0] = '005009992'
$data[130,1] = 'ALINE'
$data[130,2] = 33.41
$data[131,0] = '004435987'
$data[131,1] = 'MARCO'
$data[131,2] = 32.34
$data[132,0] = '004211911'
$data[132,1] = 'ZENILDA'
$data[132,2] = 31.9
$data[133,0] = '004207374'
$data[133,1] = 'ANGELICA'
$data[133,2] = 31.25
$data[134,0] = '005055239'
$data[134,1] = 'NORMAN'
$data[134,2] = 31.01
$data[135,0] = '005070742'
$data[135,1] = 'JUSCELINO'
$data[135,2] = 30.07
$data[136,0] = '005927101'
$data[136,1] = 'SIMONE'
$data[136,2] = 30
$data[137,0] = '004332783'
$data[137,1] = 'IRON'
$data[137,2] = 29.72
$data[138,0] = '004377415'
$data[138,1] = 'ANGELA'
$data[138,2] = 28.73
$data[139,0] = '004230529'
$data[139,1] = 'LAIS'
$data[139,2] = 28.18
$data[140,0] = '004813134'
$data[140,1] = 'MONICA'
$data[140,2] = 28.16
$data[141,0] = '005305965'
$data[141,1] = 'SIDMAR'
$data[141,2] = 28.12
$data[142,0] = '005018038'
$data[142,1] = 'ELAINE'
$data[142,2] = 27.47
$data[143,0] = '004404724'
$data[143,1] = 'LEANDRO'
$data[143,2] = 26.71
$data[144,0] = '005616259'
$data[144,1] = 'MARIA'
$data[144,2] = 26.65
$data[145,0] = '005044389'
$data[145,1] = 'CLAUDIA'
$data[145,2] = 25.45
$data[146,0] = '004350197'
$data[146,1] = 'GISELA'
$data[146,2] = 25.08
$data[147,0] = '004472760'
$data[147,1] = 'SANDRA'
$data[147,2] = 24.96
$data[148,0] = '005715733'
$data[148,1] = 'ADRIANA'
$data[148,2] = 24.5
$data[149,0] = '004756968'
$data[149,1] = 'DANIELY'
$data[149,2] = 24.09
$data[150,0] = '004243043'
$data[150,1] = 'SUELI'
$data[150,2] = 23.36
$data[151,0] = '004398174'
$data[151,1] = 'DANIELE'
$data[151,2] = 22.77
$data[152,0] = '004371857'
$data[152,1] = 'NAZARETH'
$data[152,2] = 21.52
$data[153,0] = '004388077'
$data[153,1] = 'WLADMIR'
$data[153,2] = 20.89
$data[154,0] = '004214604'
$data[154,1] = 'MARIA'
$data[154,2] = 20.75
$data[155,0] = '004467884'
$data[155,1] = 'ANA'
$data[155,2] = 20.69
$data[156,0] = '005143579'
$data[156,1] = 'GABRIEL'
$data[156,2] = 19.2
$data[157,0] = '004204255'
$data[157,1] = 'AMADO'
$data[157,2] = 18.77
$data[158,0] = '004214592'
$data[158,1] = 'MERG'
$data[158,2] = 18.64
$data[159,0] = '004920447'
$data[159,1] = 'MARILIA'
$data[159,2] = 18.19
$data[160,0] = '008032597'
$data[160,1] = 'ALESSANDRO'
$data[160,2] = 17.86
$data[161,0] = '008037529'
$data[161,1] = 'MELISSA'
$data[161,2] = 17.85
$data[162,0] = '004994036'
$data[162,1] = 'BALTASAR'
$data[162,2] = 17.62
$data[163,0] = '003497496'
$data[163,1] = 'ELISANDRA'
$data[163,2] = 17.15
$data[164,0] = '002894447'
$data[164,1] = 'JOAO'
$data[164,2] = 16.48
$data[165,0] = '005293480'
$data[165,1] = 'WAGNER'
$data[165,2] = 16.04
$data[166,0] = '004547722'
$data[166,1] = 'MARCIA'
$data[166,2] = 16.02
$data[167,0] = '004340984'
$data[167,1] = 'RENATA'
$data[167,2] = 15.62
$data[168,0] = '004422594'
$data[168,1] = 'WANDIR'
$data[168,2] = 14.67
$data[169,0] = '004455356'
$data[169,1] = 'MARCELO'
$data[169,2] = 14.59
$data[170,0] = '004570632'
$data[170,1] = 'FABRICIO'
$data[170,2] = 14.49
$data[171,0] = '005878792'
$data[171,1] = 'JUNIO'
$data[171,2] = 14.49
$data[172,0] = '005268516'
$data[172,1] = 'LUIS'
$data[172,2] = 14.27
$data[173,0] = '004526450'
$data[173,1] = 'MSD'
$data[173,2] = 14.13
$data[174,0] = '004454365'
$data[174,1] = 'RAFAEL'
$data[174,2] = 13.38
$data[175,0] = '005685089'
$data[175,1] = 'CARNEIRO'
$data[175,2] = 13.22
$data[176,0] = '005009922'
$data[176,1] = 'ANA'
$data[176,2] = 12.84
$data[177,0] = '004264780'
$data[177,1] = 'MARCELO'
$data[177,2] = 12.67
$data[178,0] = '004374943'
$data[178,1] = 'LEONARDO'
$data[178,2] = 12.08
$data[179,0] = '004551472'
$data[179,1] = 'DIEGO'
$data[179,2] = 11.88
$data[180,0] = '004472076'
$data[180,1] = 'RUBENS'
$data[180,2] = 11.8
$data[181,0] = '005133039'
$data[181,1] = 'PAULO'
$data[181,2] = 11.6
$data[182,0] = '004335031'
$data[182,1] = 'EDMUNDO'
$data[182,2] = 11.17
$data[183,0] = '005374916'
$data[183,1] = 'MARCO'
$data[183,2] = 10.89
$data[184,0] = '004216298'
$data[184,1] = 'FLORDELIZ'
$data[184,2] = 9.8
$data[185,0] = '008013889'
$data[185,1] = 'CAROLINA'
$data[185,2] = 9.72
$data[186,0] = '002828327'
$data[186,1] = 'RENAN'
$data[186,2] = 9.64
$data[187,0] = '004751154'
$data[187,1] = 'CATARINE'
$data[187,2] = 9.12
$data[188,0] = '005324981'
$data[188,1] = 'JO'
$data[188,2] = 7.69
$data[189,0] = '004308815'
$data[189,1] = 'ZELI'
$data[189,2] = 7.54
$data[190,0] = '004505474'
$data[190,1] = 'RICARDO'
$data[190,2] = 7.46
$data[191,0] = '004530494'
$data[191,1] = 'ROSANGELA'
$data[191,2] = 6.99
$data[192,0] = '004752494'
$data[192,1] = 'SERGIO'
$data[192,2] = 6.71
$data[193,0] = '008007759'
$data[193,1] = 'CRISTINA'
$data[193,2] = 6.46
$data[194,0] = '004907688'
$data[194,1] = 'HEITOR'
$data[194,2] = 6.12
$data[195,0] = '004224405'
$data[195,1] = 'MILA'
$data[195,2] = 5.88
$data[196,0] = '004228456'
$data[196,1] = 'FLASH'
$data[196,2] = 5.86
$data[197,0] = '004448501'
$data[197,1] = 'JOAO'
$data[197,2] = 5.55
$data[198,0] = '004754056'
$data[198,1] = 'BRUNO'
$data[198,2] = 5.18
$data[199,0] = '004805269'
$data[199,1] = 'CLISIA'
$data[199,2] = 5.18
$data[200,0] = '005135281'
$data[200,1] = 'RAFAEL'
$data[200,2] = 4.64
$data[201,0] = '000834301'
$data[201,1] = 'MARCUS'
$data[201,2] = 4.4
$data[202,0] = '004539779'
$data[202,1] = 'RICARDO'
$data[202,2] = 4.37
$data[203,0] = '008012870'
$data[203,1] = 'ANA'
$data[203,2] = 4.3
$data[204,0] = '004165515'
$data[204,1] = 'MAURO'
$data[204,2] = 4.2
$data[205,0] = '008004995'
$data[205,1] = 'JOSE'
$data[205,2] = 3.74
$data[206,0] = '004328934'
$data[206,1] = 'VALERIA'
$data[206,2] = 3.55
$data[207,0] = '004352384'
$data[207,1] = 'BRASFORT'
$data[207,2] = 3.13
$data[208,0] = '004181486'
$data[208,1] = 'ANDREA'
$data[208,2] = 3.11
$data[209,0] = '005170415'
$data[209,1] = 'MONICA'
$data[209,2] = 2.84
$data[210,0] = '005140667'
$data[210,1] = 'MATEUS'
$data[210,2] = 2.2
$data[211,0] = '005022526'
$data[211,1] = 'ALEXANDRE'
$data[211,2] = 1.7
$data[212,0] = '000431814'
$data[212,1] = 'GUILHERME'
$data[212,2] = 1.1
$data[213,0] = '004360430'
$data[213,1] = 'VIOMAR'
$data[213,2] = 1
$data[214,0] = '004486497'
$data[214,1] = 'ELENA'
$data[214,2] = 0.96
$data[215,0] = '004115403'
$data[215,1] = 'HEBERT'
$data[215,2] = 0.88
$data[216,0] = '005660155'
$data[216,1] = 'CAROLINA'
$data[216,2] = 0.85
$data[217,0] = '004223502'
$data[217,1] = 'BRUNA'
$data[217,2] = 0.78
$data[218,0] = '002687737'
$data[218,1] = 'JOSE'
$data[218,2] = 0.71
$data[219,0] = '004587511'
$data[219,1] = 'CARLOS'
$data[219,2] = 0.69
$data[220,0] = '004380948'
$data[220,1] = 'LUISA'
$data[220,2] = 0.62
$data[221,0] = '004473942'
$data[221,1] = 'DAIANNE'
$data[221,2] = 0.62
$data[222,0] = '003894173'
$data[222,1] = 'ANDREA'
$data[222,2] = 0.48
$data[223,0] = '004453302'
$data[223,1] = 'ISABELLA'
$data[223,2] = 0.39
$data[224,0] = '004638738'
$data[224,1] = 'GABRIEL'
$data[224,2] = 0.29
$data[225,0] = '004278033'
$data[225,1] = 'DAISY'
$data[225,2] = 0.21
$data[226,0] = '005662526'
$data[226,1] = 'AGUINALDO'
$data[226,2] = 0.18
$data[227,0] = '004432455'
$data[227,1] = 'LUCIANA'
$data[227,2] = 0.17
$data[228,0] = '002694089'
$data[228,1] = 'VITOR'
$data[228,2] = 0.16
$data[229,0] = '004357159'
$data[229,1] = 'JOAO'
$data[229,2] = 0.15
$data[230,0] = '004320840'
$data[230,1] = 'NATALIA'
$data[230,2] = 0.14
$data[231,0] = '001000288'
$data[231,1] = 'ISABELLA'
$data[231,2] = 0.13
$data[232,0] = '005530256'
$data[232,1] = 'CAROLINA'
$data[232,2] = 0.1
$data[233,0] = '004451996'
$data[233,1] = 'ADRIANO'
$data[233,2] = 0.09
$data[234,0] = '005047946'
$data[234,1] = 'GABRIEL'
$data[234,2] = 0.09
$data[235,0] = '004223226'
$data[235,1] = 'YESHUA'
$data[235,2] = 0.04
$data[236,0] = '004281300'
$data[236,1] = 'FRANKLIN'
$data[236,2] = 0.04
$data[237,0] = '005274028'
$data[237,1] = 'RAFAEL'
$data[237,2] = 0.04
$data[238,0] = '004329229'
$data[238,1] = 'GABRIEL'
$data[238,2] = 0.03
$data[239,0] = '004213373'
$data[239,1] = 'ALEXANDRE'
$data[239,2] = 0.02
$data[240,0] = '004339183'
$data[240,1] = 'JALISON'
$data[240,2] = 0.02
$data[241,0] = '004870976'
$data[241,1] = 'HFR'
$data[241,2] = 0.02
$data[242,0] = '000938440'
$data[242,1] = 'BASE'
$data[242,2] = 0.01
$data[243,0] = '002878817'
$data[243,1] = 'GUILHERME'
$data[243,2] = 0.01
$data[244,0] = '004400000'
$data[244,1] = 'VILMA'
$data[244,2] = 0.01
$data[245,0] = '004976625'
$data[245,1] = 'NORTON'
$data[245,2] = 0.01
$data[246,0] = '005105970'
$data[246,1] = 'VERA'
$data[246,2] = -11747.84

$ws.Range("A2:C248").Value = $data

# Row after the data (blank separator row) - ensure it has no leftover content
$ws.Range("A249:C249").ClearContents()

# Footer row with the applied-filters description (single cell, column A only)
$ws.Range("A250").Value = "Filtros aplicados:`nDataFim é (Em branco)`nnr_saldo_disponivel não é 0`nPosição é Posição D-1`nCARTEIRA não está em branco`nDataFim é (Em branco)`nNR_CONTA não está em branco`nTIPO_LANCAMENTO não é ED, ET ou Liquidação Doador"

Write-Output "done"
